# Apply the amova.docx edits:
#  1. Widen the first table column (4777 -> 5045 twips, i.e. 238.85pt -> 252.25pt)
#  2. Rename "Among populations within urban/rural groups" -> "Among sampling sites within urban/rural groups"
#  3. Bump the "Within populations" row height (612 -> 614 twips, i.e. 30.6pt -> 30.7pt)
#  4. Rename "Within populations" -> "Within sampling sites"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. First column width: 5045 twips = 252.25 points (gridCol w:w is in twentieths of a point)
$t.Columns.Item(1).Width = 252.25

# 2 & 4. Text replacements
$d.Content.Find.Execute("Among populations within urban/rural groups", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Among sampling sites within urban/rural groups", 2)
$d.Content.Find.Execute("Within populations", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Within sampling sites", 2)

# 3. Row height: 614 twips = 30.7 points, on the row that now reads "Within sampling sites"
$t.Rows.Item(4).Height = 30.7
